# "Menambah flow panel 1 dan panel 2" - add ONT-failure / no-GPON / no-port
# flow rows (panel 1 & panel 2) highlighted in yellow, plus ODP labels for a
# few existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: move the ODP label from H11 to G11 ("ONT Gagal") -------------
$ws.Range("G11").Value = "ONT Gagal"
$ws.Range("H11").ClearContents()

# --- Row 12: new flow entry (panel 1), highlighted yellow like row 6 ------
$ws.Range("B6:F6").Copy() | Out-Null
$ws.Range("B12:F12").PasteSpecial(-4122) | Out-Null

$ws.Range("G12").Value = 23045811
$ws.Range("G12").Interior.Color = 65535
$ws.Range("G12").VerticalAlignment = -4160

$ws.Range("H12").Value = "ODP-SUD-FCL/08"
$ws.Range("H12").Interior.Color = 65535
$ws.Range("H12").VerticalAlignment = -4160

$ws.Range("I12").Value = 16
$ws.Range("I12").Interior.Color = 65535
$ws.Range("I12").VerticalAlignment = -4160

$ws.Range("J12").Value = 3
$ws.Range("J12").Interior.Color = 65535
$ws.Range("J12").VerticalAlignment = -4160

# --- Row 13: ODP label + shrink the row back to normal height -------------
$ws.Range("G13").Value = "GPON Tidak Ada"
$ws.Rows.Item(13).RowHeight = 15

# --- Row 14: ODP label + max/ready port, shrink row height ----------------
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = "ODP-SUD-FBC/074"
$ws.Range("I14").Value = 8
$ws.Range("J14").Value = 4
$ws.Rows.Item(14).RowHeight = 15

# --- Row 15: ODP label -----------------------------------------------------
$ws.Range("G15").Value = "Port Tidak Ada"

# --- Row 16: ODP label + max/ready port -----------------------------------
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$ws.Range("H16").Value = "ODP-SUD-FAB/041"
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 4

# --- Row 17-20: ODP labels --------------------------------------------------
$ws.Range("G17").Value = "GPON Tidak Ada"
$ws.Range("G18").Value = "Port Tidak Ada"
$ws.Range("G19").Value = "Port Tidak Ada"
$ws.Range("G20").Value = "Port Tidak Ada"

# --- Row 21: ODP label + max/ready port -----------------------------------
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H21").PasteSpecial(-4122) | Out-Null
$ws.Range("H21").Value = "ODP-SUD-FAZ/041"
$ws.Range("I21").Value = 8
$ws.Range("J21").Value = 2

# --- Window / selection state ----------------------------------------------
$ws.Range("I13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
